# Apply KPI computation changes across Productdata, Capacity and ProcessingTime sheets.
$wb = $excel.ActiveWorkbook

# --- Productdata sheet: update StartingInventories (C) and SetupCosts (E) ---
$wsProduct = $wb.Worksheets.Item("Productdata")

$productC = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 0
    6 = 0
    7 = 40834
    8 = 19456
    9 = 58601
    10 = 7482
    11 = 0
    12 = 0
    13 = 0
    14 = 0
}

$productE = @{
    2 = 1008.888888888889
    3 = 205.6666666666666
    4 = 102.7777777777778
    5 = 410.4444444444443
    6 = 446.3583333333332
    7 = 205.6666666666666
    8 = 102.7777777777778
    9 = 820.8888888888887
    10 = 35.91388888888888
    11 = 179
    12 = 89.44444444444443
    13 = 357.111111111111
    14 = 388.3583333333333
}

foreach ($row in 2..14) {
    $wsProduct.Range("C$row").Value = $productC[$row]
    $wsProduct.Range("E$row").Value = $productE[$row]
}

# --- Capacity sheet: update capacity values (B) ---
$wsCapacity = $wb.Worksheets.Item("Capacity")

$capacityB = @{
    2 = 908000
    3 = 80000
    4 = 100000
    5 = 400000
    6 = 348000
    7 = 80000
    8 = 80000
    9 = 160000
    10 = 7000
    11 = 120000
    12 = 60000
    13 = 80000
    14 = 261000
}

foreach ($row in 2..14) {
    $wsCapacity.Range("B$row").Value = $capacityB[$row]
}

# --- ProcessingTime sheet: update diagonal processing times ---
$wsProcessing = $wb.Worksheets.Item("ProcessingTime")

$processingDiag = @{
    2  = @{ Col = "B"; Value = 4 }
    3  = @{ Col = "C"; Value = 2 }
    4  = @{ Col = "D"; Value = 5 }
    5  = @{ Col = "E"; Value = 5 }
    7  = @{ Col = "G"; Value = 2 }
    8  = @{ Col = "H"; Value = 4 }
    9  = @{ Col = "I"; Value = 2 }
    10 = @{ Col = "J"; Value = 1 }
    11 = @{ Col = "K"; Value = 3 }
    13 = @{ Col = "M"; Value = 1 }
    14 = @{ Col = "N"; Value = 3 }
}

foreach ($row in $processingDiag.Keys) {
    $entry = $processingDiag[$row]
    $addr = "$($entry.Col)$row"
    $wsProcessing.Range($addr).Value = $entry.Value
}
